# ConfigSheet.xlsx - "credit request customer create 14/12/2017"
#
# The test fixture's "Sheet" tab holds key/value pairs that drive Selenium
# API tests. The "AddNewPartnerOrganisationTest" block (domainName /
# tradeLicenseNo / vatId) and the "CustomerTest" block (customerdomainName /
# companyName) need fresh, previously-unused identifiers so a re-run of the
# credit-request / customer-create flow doesn't collide with entities that
# already exist server-side.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")
$ws.Activate()

# AddNewPartnerOrganisationTest block
$ws.Range("B15").Value = "burhani7003"      # domainName
$ws.Range("B16").Value = "tno400002340"     # tradeLicenseNo
$ws.Range("B17").Value = "vat390002460"     # vatId

# CustomerTest block
$ws.Range("B20").Value = "customer7014"     # customerdomainName
$ws.Range("B21").Value = "CompanyName7003"  # companyName

# Matches the author's final cursor position recorded in the saved view state.
$ws.Range("B17").Select()
